# Add a new row to the Table1 listing ("hebrew_text" / "hebrew_option" / "label")
# for the P_ADHD_Restless question variant that is missing a comma, then leave
# the selection where the user last clicked (A10), matching the saved session
# state recorded in the workbook after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# The table is sorted by the "label" column; the new row sorts to sheet row 3
# (the 2nd data row), right after the existing "P_ADHD_Distracted" row and
# before the existing (comma-variant) "P_ADHD_Restless" row. Insert a fresh
# worksheet row there so every following row shifts down by one, matching the
# table's existing rows untouched.
$ws.Rows.Item(3).Insert() | Out-Null

$ws.Cells.Item(3, 1).Value2 = "למיטב ידיעתי - ברגע זה הילד/ה שלי מרגיש/ה חסר/ת מנוחה"
$ws.Cells.Item(3, 3).Value2 = "P_ADHD_Restless"

# Grow the table definition (and its AutoFilter) to cover the new row.
$lo.Resize($ws.Range("A1:C35")) | Out-Null

# Restore the cursor/selection position recorded in the saved workbook.
$ws.Range("A10").Select() | Out-Null
